# Add two new columns, I ("I0") and J ("IF"), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers: copy the existing header formatting (bold font + border + centered
# alignment, same style used by B1:H1) from H1 onto I1:J1, then set the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-39: I = dS0-ish "I0" metric, J = "IF" metric (plain numbers,
# same unstyled look as the rest of the numeric columns).
$values = @(
  @(6,6), @(7,7), @(6,6), @(6,6), @(8,8), @(8,8), @(7,7), @(11,11), @(7,7),
  @(6,6), @(8,8), @(7,7), @(7,7), @(7,7), @(5,6), @(8,8), @(7,7), @(11,11),
  @(7,7), @(9,9), @(7,7), @(8,8), @(7,7), @(7,7), @(9,10), @(7,7), @(7,7),
  @(7,7), @(6,7), @(7,7), @(9,9), @(8,8), @(4,4), @(5,6), @(7,7), @(9,9),
  @(3,3), @(8,8)
)

for ($k = 0; $k -lt $values.Length; $k++) {
  $row = $k + 2
  $pair = $values[$k]
  $ws.Cells.Item($row, 9).Value = $pair[0]
  $ws.Cells.Item($row, 10).Value = $pair[1]
}
